# Add a new "2020" column (column Q) to the SDG 1.5.1 table, mirroring the
# formatting of the existing 2019 column (column P) and filling in the new
# year's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (number format, font, borders, alignment, etc.) from the
# 2019 column (P, rows 3-34: divider row through the last data row) onto the
# new 2020 column (Q) before writing the new values.
$ws.Range("P3:P34").Copy($ws.Range("Q3:Q34"))

# Header
$ws.Range("Q4").Value2 = 2020

# Values for the new 2020 column, row by row. "-" means "no data" (same
# convention used throughout the rest of the table).
$q = @{
    5  = 51
    6  = 29
    7  = 22
    8  = 5
    9  = 3
    10 = 2
    11 = 15
    12 = 9
    13 = 5
    14 = "-"
    15 = "-"
    16 = "-"
    17 = "-"
    18 = "-"
    19 = "-"
    20 = 7
    21 = 7
    22 = "-"
    23 = "-"
    24 = "-"
    25 = "-"
    26 = 24
    27 = 10
    28 = 14
    29 = "-"
    30 = "-"
    31 = "-"
    32 = "-"
    33 = "-"
    34 = "-"
}

foreach ($row in $q.Keys) {
    $ws.Cells.Item($row, 17).Value2 = $q[$row]
}

# Update the active selection to match the author's final cursor position.
$ws.Range("Q35").Select() | Out-Null
